# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G ("K") holds per-game strikeout counts. The values below replace the
# previously stored "Strike#" figures with the actual strikeout totals (K).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 5
    3  = 2
    4  = 7
    5  = 7
    6  = 5
    7  = 3
    8  = 2
    9  = 4
    10 = 9
    11 = 6
    12 = 5
    13 = 5
    14 = 4
    15 = 6
    16 = 4
    17 = 5
    18 = 4
    19 = 5
    20 = 4
    21 = 10
    22 = 3
    23 = 11
    24 = 7
    25 = 5
    26 = 5
    27 = 8
    28 = 4
    29 = 8
    30 = 4
    31 = 6
    32 = 4
    33 = 2
    34 = 4
    35 = 6
    36 = 5
    37 = 6
    38 = 5
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
